$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Top block (A2) drives A3:A5 via existing formulas; just bump the input ---
$ws.Range("A2").Value = 60

# --- Row 17 ---
$ws.Range("B17").Value = 5

# --- Row 19: new B value + relabelled C (string) ---
$ws.Range("B19").Value = 3
$ws.Range("C19").Value = "wateringRecordIndex"

# --- Row 20: new B value + relabelled C (string) ---
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = "backlightMode"

# --- Row 21: new B value (C stays "wateringSettings (3 times)") ---
$ws.Range("B21").Value = 30
$ws.Range("C21").Value = "wateringSettings (3 times)"

# --- Row 22: new B value (C stays "wateringStatus") ---
$ws.Range("B22").Value = 13
$ws.Range("C22").Value = "wateringStatus"

# --- Row 23: new B value + relabelled C (string) ---
$ws.Range("B23").Value = 2
$ws.Range("C23").Value = "cumulativeRunningHour"

# --- Row 24: now gets B and C values for the first time ---
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = "hourOfDay"

# --- New row 28, continuing the running-total formula series ---
$ws.Range("A28").Formula = "=A27+B27"
